# Guiglia.xlsx update: a new day (2020-10-28 / serial 44235) of data was
# inserted into the time series, and a new day (serial 44257) was appended
# at the end. The 7-day rolling-sum columns (C, D) were recomputed for all
# rows whose trailing 7-day window is affected by the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row for serial date 44235 between rows 92 and 93 ---
# This shifts the existing rows 93:113 down to 94:114.
$ws.Rows.Item(93).Insert()

# The freshly inserted row doesn't inherit the date column's cell style,
# so copy formatting (incl. the date number format) from the row below.
$ws.Range("A94").Copy()
$ws.Range("A93").PasteSpecial(-4122)

# Populate the newly inserted row (93) with its data.
$ws.Range("A93").Value = 44235
$ws.Range("B93").Value = 4
$ws.Range("C93").Value = 22
$ws.Range("D93").Value = 559.5116988809766

# --- 2. Recompute the rolling-sum columns (C, D) for rows whose 7-day ---
# --- window now includes the newly inserted day.                     ---
$ws.Range("C90").Value = 30
$ws.Range("D90").Value = 762.970498474059

$ws.Range("C91").Value = 21
$ws.Range("D91").Value = 534.0793489318413

$ws.Range("C92").Value = 22
$ws.Range("D92").Value = 559.5116988809766

$ws.Range("C94").Value = 17
$ws.Range("D94").Value = 432.3499491353002

$ws.Range("C95").Value = 14
$ws.Range("D95").Value = 356.0528992878942

# Row that shifted from the old row 111 (now row 112) previously had no
# rolling-sum figures yet; they can now be computed.
$ws.Range("C112").Value = 29
$ws.Range("D112").Value = 737.5381485249237

# --- 3. Append a new row for serial date 44257 at the end of the table ---
# (after the insert above, the last existing data row is now row 114,
# corresponding to the original last row / serial date 44256)
$ws.Range("A114").Copy()
$ws.Range("A115").PasteSpecial(-4122)
$ws.Range("A115").Value = 44257
$ws.Range("B115").Value = 3
# C/D have no rolling-sum value yet for this brand-new day (same as the
# two rows before it), so leave them as empty placeholders.
$ws.Range("C115").Value = ""
$ws.Range("D115").Value = ""
